$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-12 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-13 Thursday", 2) | Out-Null
$d.Content.Find.Execute("663÷8=82, 7", $true, $false, $false, $false, $false, $true, 1, $false, "280÷6=46, 4", 2) | Out-Null
$d.Content.Find.Execute("669÷8=83, 5", $true, $false, $false, $false, $false, $true, 1, $false, "996÷4=249, 0", 2) | Out-Null
$d.Content.Find.Execute("409÷8=51, 1", $true, $false, $false, $false, $false, $true, 1, $false, "532÷2=266, 0", 2) | Out-Null
$d.Content.Find.Execute("559÷2=279, 1", $true, $false, $false, $false, $false, $true, 1, $false, "753÷7=107, 4", 2) | Out-Null
$d.Content.Find.Execute("347÷4=86, 3", $true, $false, $false, $false, $false, $true, 1, $false, "930÷3=310, 0", 2) | Out-Null
$d.Content.Find.Execute("148÷2=74, 0", $true, $false, $false, $false, $false, $true, 1, $false, "247÷7=35, 2", 2) | Out-Null
$d.Content.Find.Execute("550÷9=61, 1", $true, $false, $false, $false, $false, $true, 1, $false, "437÷4=109, 1", 2) | Out-Null
$d.Content.Find.Execute("710÷6=118, 2", $true, $false, $false, $false, $false, $true, 1, $false, "966÷2=483, 0", 2) | Out-Null
$d.Content.Find.Execute("362÷4=90, 2", $true, $false, $false, $false, $false, $true, 1, $false, "746÷9=82, 8", 2) | Out-Null
$d.Content.Find.Execute("885÷3=295, 0", $true, $false, $false, $false, $false, $true, 1, $false, "902÷5=180, 2", 2) | Out-Null
$d.Content.Find.Execute("914÷4=228, 2", $true, $false, $false, $false, $false, $true, 1, $false, "976÷4=244, 0", 2) | Out-Null
$d.Content.Find.Execute("470÷6=78, 2", $true, $false, $false, $false, $false, $true, 1, $false, "695÷3=231, 2", 2) | Out-Null
$d.Content.Find.Execute("356÷3=118, 2", $true, $false, $false, $false, $false, $true, 1, $false, "156÷9=17, 3", 2) | Out-Null
$d.Content.Find.Execute("939÷3=313, 0", $true, $false, $false, $false, $false, $true, 1, $false, "545÷6=90, 5", 2) | Out-Null
$d.Content.Find.Execute("830÷2=415, 0", $true, $false, $false, $false, $false, $true, 1, $false, "494÷9=54, 8", 2) | Out-Null
$d.Content.Find.Execute("303÷5=60, 3", $true, $false, $false, $false, $false, $true, 1, $false, "104÷3=34, 2", 2) | Out-Null
$d.Content.Find.Execute("482÷3=160, 2", $true, $false, $false, $false, $false, $true, 1, $false, "942÷9=104, 6", 2) | Out-Null
$d.Content.Find.Execute("324÷4=81, 0", $true, $false, $false, $false, $false, $true, 1, $false, "886÷2=443, 0", 2) | Out-Null
$d.Content.Find.Execute("574÷8=71, 6", $true, $false, $false, $false, $false, $true, 1, $false, "474÷6=79, 0", 2) | Out-Null
$d.Content.Find.Execute("316÷6=52, 4", $true, $false, $false, $false, $false, $true, 1, $false, "717÷5=143, 2", 2) | Out-Null
$d.Content.Find.Execute("812÷8=101, 4", $true, $false, $false, $false, $false, $true, 1, $false, "306÷4=76, 2", 2) | Out-Null
$d.Content.Find.Execute("405÷3=135, 0", $true, $false, $false, $false, $false, $true, 1, $false, "196÷7=28, 0", 2) | Out-Null
$d.Content.Find.Execute("102÷9=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "561÷6=93, 3", 2) | Out-Null
$d.Content.Find.Execute("687÷5=137, 2", $true, $false, $false, $false, $false, $true, 1, $false, "113÷5=22, 3", 2) | Out-Null
$d.Content.Find.Execute("879÷3=293, 0", $true, $false, $false, $false, $false, $true, 1, $false, "560÷6=93, 2", 2) | Out-Null
